$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6982125043869019
$ws.Range("B1").Value = 1.913483142852783
$ws.Range("C1").Value = 5.700013160705566
$ws.Range("D1").Value = 1.534625172615051
$ws.Range("E1").Value = 0.6905362010002136
